$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "45÷4="
$t.Cell(1, 2).Range.Text = "17÷2="
$t.Cell(1, 3).Range.Text = "38÷5="
$t.Cell(1, 4).Range.Text = "54÷5="
$t.Cell(1, 5).Range.Text = "70÷2="

# Row 5
$t.Cell(5, 1).Range.Text = "11÷4="
$t.Cell(5, 2).Range.Text = "24÷5="
$t.Cell(5, 3).Range.Text = "73÷4="
$t.Cell(5, 4).Range.Text = "94÷6="
$t.Cell(5, 5).Range.Text = "28÷4="

# Row 9
$t.Cell(9, 1).Range.Text = "47÷8="
$t.Cell(9, 2).Range.Text = "37÷5="
$t.Cell(9, 3).Range.Text = "25÷2="
$t.Cell(9, 4).Range.Text = "39÷9="
$t.Cell(9, 5).Range.Text = "68÷3="

# Row 13
$t.Cell(13, 1).Range.Text = "43÷4="
$t.Cell(13, 2).Range.Text = "64÷3="
$t.Cell(13, 3).Range.Text = "35÷6="
$t.Cell(13, 4).Range.Text = "12÷4="
$t.Cell(13, 5).Range.Text = "91÷5="

# Row 17
$t.Cell(17, 1).Range.Text = "80÷7="
$t.Cell(17, 2).Range.Text = "94÷6="
$t.Cell(17, 3).Range.Text = "65÷9="
$t.Cell(17, 4).Range.Text = "31÷2="
$t.Cell(17, 5).Range.Text = "82÷3="
